$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.775.69"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "3.804.83"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "3.802.45"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "4.444.03"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "3.760.03"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.58%  "
$ws.Range("D18").Value = "67.779.46"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "461.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("E22").Value = "  -7.69%  "
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "3.951.70"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +4.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "391.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -4.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.30%  "
